$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3: H3 -96 -> -97, I3 "LAST UPDATE" text 03-Nov-2025 -> 04-Nov-2025
$ws.Range("H3").Value = -97
# I3 holds the date as literal text (not a real date value), so force the
# cell to Text format first - otherwise Excel auto-converts the date-like
# string into a date serial number, which changes both the value and type.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "04-Nov-2025"

# Row 4: H4 609 -> 608, I4 "LAST UPDATE" text 03-Nov-2025 -> 04-Nov-2025
$ws.Range("H4").Value = 608
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "04-Nov-2025"
